# Insert a new weekly price-report row for Cilantro (Terminal La Palmera de
# La Serena) at row 125, pushing the existing rows 125-159 down to 126-160.
# The new row carries the same Min/Max/Avg/Price-per-kg figures as the
# (former) row 125, but with an updated date and a new Volumen value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 125..159 down to 126..160, opening up a blank row 125.
$ws.Rows(125).Insert()

# Populate the newly inserted row 125.
$ws.Range("A125").Value = 8
$ws.Range("B125").Value = "Terminal La Palmera de La Serena"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44782
$ws.Range("E125").Value = 4
$ws.Range("F125").Value = 100112040
$ws.Range("G125").Value = "Cilantro"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 2400
$ws.Range("K125").Value = 2000
$ws.Range("L125").Value = 2500
$ws.Range("M125").Value = 2250
$ws.Range("N125").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O125").Value = "Provincia del Elquí"
$ws.Range("P125").Value = 1500
$ws.Range("Q125").Value = 1.5
$ws.Range("R125").Value = "Hortaliza"
